$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.099.43'
$ws.Range("E2").Value = '  -0.85%  '

$ws.Range("D3").Value = '3.397.28'
$ws.Range("E3").Value = '  -3.71%  '

$ws.Range("E4").Value = '  +0.06%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '579.53'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.37%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '136.16'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -5.26%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = '3.397.45'
$ws.Range("E8").Value = '  -3.64%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.492'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.28%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.11'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -8.91%  '

$ws.Range("E11").Value = '  -11.54%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.370'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -8.66%  '

$ws.Range("D13").Value = '3.978.63'
$ws.Range("E13").Value = '  -3.68%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.0000176'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -11.65%  '

$ws.Range("E15").Value = '  -1.72%  '

$ws.Range("D16").Value = '3.414.57'
$ws.Range("E16").Value = '  -3.04%  '

$ws.Range("D17").Value = '65.105.77'
$ws.Range("E17").Value = '  -0.87%  '

$ws.Range("E18").Value = '  -10.38%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '9.55'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -12.58%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '5.81'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -6.44%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '13.45'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -6.13%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '378.93'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -8.92%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.547'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -8.80%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.11%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '71.78'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -7.11%  '

$ws.Range("D26").Value = '3.532.50'
$ws.Range("E26").Value = '  -3.75%  '

$ws.Range("E27").Value = '  -10.15%  '

$ws.Range("E28").Value = '  +0.09%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '6.98'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -9.87%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.19'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -10.39%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.96'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -10.88%  '

$ws.Range("D32").Value = '3.409.05'
$ws.Range("E32").Value = '  -3.39%  '

$ws.Range("E33").Value = '  -0.03%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.141'
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '22.71'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -6.80%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '169.40'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -2.85%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '6.64'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -12.24%  '

$ws.Range("E38").Value = '  -11.57%  '

$ws.Range("E39").Value = '  -8.03%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '4.64'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -12.47%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.0750'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -8.83%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.803'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -6.52%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '43.23'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -4.53%  '

$ws.Range("E44").Value = '  +0.17%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '4.32'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -15.20%  '

$ws.Range("E46").Value = '  -10.69%  '

$ws.Range("E47").Value = '  -0.35%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '21.92'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.93%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '6.42'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -9.51%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.04'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -14.37%  '

$ws.Range("D51").Value = '2.161.47'
$ws.Range("E51").Value = '  -8.13%  '
